$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Full Name"
$ws.Range("E1").Value = "Email"
$ws.Range("E2").Select()
